# Apply cryptos.xlsx data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.311.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.930.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7490"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3184"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08049"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.82%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.950.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.405"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.321.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "252.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007923"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.139.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.662"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.589"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1300"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.188"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.560"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.427"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.144"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05257"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.320"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7590"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.775"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01956"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.798"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.517"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4537"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.975"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8442"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.711"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.78%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.105.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.63%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1219"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.86%  "
